# Actualizacion en dias de nomina y gral
# Varios archivos para los dias de nomina e impresiones en prima de vacaciones

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PagoPlanta20141001")

# --- D44: payroll account number was missing (0), now populated ---
$ws.Cells.Item(44, 4).Value = 6100909057

# --- B54 / B55: updated payroll day amounts ---
$ws.Cells.Item(54, 2).Value = 2440738
$ws.Cells.Item(55, 2).Value = 1729103

# --- Rows 69/70: the two employee records (LOPEZ GONZALEZ ORLANDO and
#     LOPEZ GONZALEZ YELENIS YANIT) swap places completely ---
$ws.Cells.Item(69, 1).Value = 17805992
$ws.Cells.Item(69, 2).Value = 2452277
$ws.Cells.Item(69, 3).Value = "LOPEZ GONZALEZ ORLANDO "
$ws.Cells.Item(69, 4).Value = 236070213244

$ws.Cells.Item(70, 1).Value = 40916440
$ws.Cells.Item(70, 2).Value = 2730202
$ws.Cells.Item(70, 3).Value = "LOPEZ GONZALEZ YELENIS YANIT"
$ws.Cells.Item(70, 4).Value = 236070069976

# --- Rows 93/94: OCHOA REDONDO ANA JOSEFA and OCHOA REDONDO YAMILKA
#     KARINA swap row order (their IDs follow them) and get new
#     updated payroll amounts in column B ---
$ws.Cells.Item(93, 1).Value = 40918307
$ws.Cells.Item(93, 2).Value = 2298409
$ws.Cells.Item(93, 3).Value = "OCHOA REDONDO ANA JOSEFA"
$ws.Cells.Item(93, 4).Value = 236070228622

$ws.Cells.Item(94, 1).Value = 40929923
$ws.Cells.Item(94, 2).Value = 2356110
$ws.Cells.Item(94, 3).Value = "OCHOA REDONDO YAMILKA KARINA"
$ws.Cells.Item(94, 4).Value = 236070218003
